# issue #5: stock data output to json file
#
# The "股票" (stock) worksheet (sheet 4) gains a new "property_category"
# column, inserted right before the existing "date" column. Every data
# row on that sheet gets the literal value "stock" in the new column,
# while the previously adjacent "date" / "legislator_name" /
# "legislator_id" columns shift one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# Insert a new column H; this shifts the old H (date), I (legislator_name)
# and J (legislator_id) columns to I, J, K respectively and keeps their
# values intact.
$ws.Columns("H:H").Insert()

# New header for the inserted column.
$ws.Range("H1").Value = "property_category"

# Fill in the new column's value for every existing data row (rows 2-7
# hold the six stock holdings).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = "stock"
}
